$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestResults")

$ws.Range("A1").Value = "Test Case ID"
$ws.Range("B1").Value = "Test scenario"
$ws.Range("C1").Value = "Status"

$ws.Range("A2").Value = "iProc_TC_ID_1"
$ws.Range("B2").Value = "@Smoke Verify Elumina Login"
$ws.Range("C2").Value = "passed"

$ws.Range("A3").Value = "iProc_TC_ID_1A"
$ws.Range("B3").Value = "@Smoke Verify Elumina Registration"
$ws.Range("C3").Value = "passed"

$ws.Range("A4").Value = "iProc_TC_ID_3"
$ws.Range("B4").Value = "@Smoke Verify CandidatesInvExam"
$ws.Range("C4").Value = "passed"

$ws.Range("A5").Value = "iProc_TC_ID_5"
$ws.Range("B5").Value = "@SmokeValidation of user authentication by valid Candidate Creadentials`n"
$ws.Range("C5").Value = "passed"

$ws.Range("A6").Value = "iProc_TC_ID_8"
$ws.Range("B6").Value = "@Smoke Validation of candidate choosing proctoring exam in dashboard"
$ws.Range("C6").Value = "passed"

$ws.Range("A7").Value = "iProc_TC_ID_23"
$ws.Range("B7").Value = "@iProctorlink Verify CandidatesExam"
$ws.Range("C7").Value = "passed"

$ws.Range("A8").Value = "iProc_TC_ID_25"
$ws.Range("B8").Value = "@iProctorlink Verify CandidatesExam"
$ws.Range("C8").Value = "passed"

$ws.Range("A9").Value = "iProc_TC_ID_28"
$ws.Range("B9").Value = "@iProctorlink Verify Elumina Invigilator Dashboard"
$ws.Range("C9").Value = "passed"

$ws.Range("A10").Value = "iProc_TC_ID_56"
$ws.Range("B10").Value = "@Smoke Verify Validation of Invigilator Dashboard Proctor"
$ws.Range("C10").Value = "passed"

$ws.Range("A11").Value = "iProc_TC_ID_57"
$ws.Range("B11").Value = "@Smoke Validation of Navigating to an exam from the dashboard to invigilate"
$ws.Range("C11").Value = "passed"

$ws.Range("A12").Value = "iProc_TC_ID_58"
$ws.Range("B12").Value = "@iProctorlink Verify Validation of `"Start Exam`" (All Candidates) Proctor "
$ws.Range("C12").Value = "timedOut"

$ws.Range("A13").Value = "iProc_TC_ID_59"
$ws.Range("B13").Value = "@iProctorlink Verify Validation of `"Lock Exam`" from Live monitor Proctor"
$ws.Range("C13").Value = "passed"

$ws.Range("A14").Value = "iProc_TC_ID_61"
$ws.Range("B14").Value = "@iProctorlink Verify Validation of `"Resume Exam`" from Live monitor Proctor"
$ws.Range("C14").Value = "passed"

$ws.Range("A15").Value = "iProc_TC_ID_64"
$ws.Range("B15").Value = "@Smoke Verify Validation of `"Mark Attendance`" (All Candidates) Proctor"
$ws.Range("C15").Value = "failed"

$ws.Range("A16").Value = "iProc_TC_ID_70"
$ws.Range("B16").Value = "@Smoke Validation of Questions answered / Inprogress on the RHS of the Candidate page"
$ws.Range("C16").Value = "passed"

$ws.Range("A17").Value = "iProc_TC_ID_71"
$ws.Range("B17").Value = "@Smoke Validation of all the events generated on the RHS of the Candidate page"
$ws.Range("C17").Value = "passed"

$ws.Range("A18").Value = "Exam_Prerequisit_ID_01"
$ws.Range("B18").Value = "@iProctorlink Verify Create Exam With Content Section and Content Section Page"
$ws.Range("C18").Value = "failed"

$ws.Range("A19").Value = "Reg_Prerequisit_ID_01A"
$ws.Range("B19").Value = "@iProctorlink Verify Elumina Registration"
$ws.Range("C19").Value = "failed"

$ws.Range("A20").Value = "iProc_TC_ID_40"
$ws.Range("B20").Value = "@iProctorlink Verify Elumina Invigilator Dashboard"
$ws.Range("C20").Value = "timedOut"

$ws.Range("A21").Value = "@Smoke Verify Elumina Login and Create Exam"
$ws.Range("C21").Value = "passed"

$ws.Range("A22").Value = "@Smoke Verify Elumina RegistrationInv and add User and Invigilator"
$ws.Range("C22").Value = "passed"

$ws.Range("A23").Value = "iProc_TC_ID_11"
$ws.Range("B23").Value = "@Smoke Validation of `"I Authorise`" checkbox - To access Webcam, Microphone & Terms & Condition"
$ws.Range("C23").Value = "passed"

$ws.Range("A24").Value = "iProc_TC_ID_34"
$ws.Range("B24").Value = "@Smoke Validation of Exam section page (Offline Exam validation)"
$ws.Range("C24").Value = "passed"

$ws.Range("A25").Value = "iProc_TC_ID_39"
$ws.Range("B25").Value = "@Smoke Validation of submitting when the Candidate has not answered all Questions"
$ws.Range("C25").Value = "passed"

